$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("F2").Value = 1725526800
$ws.Range("G2").Value = 1727686800
$ws.Range("H2").Value = "SIN MOVIL"
$ws.Range("I2").Value = ""
$ws.Range("AK2").Value = "Existe en la BD"

# Row 3
$ws.Range("F3").Value = 1725526800
$ws.Range("G3").Value = 1725526800
$ws.Range("AK3").Value = "Existe en la BD"

# Row 4
$ws.Range("F4").Value = 1725526800
$ws.Range("G4").Value = 1726131600
$ws.Range("AK4").Value = "Existe en la BD"
